# Update the cryptos price/volume table with refreshed quotes.
# Note: several "Price" column values are plain decimal-looking strings
# (e.g. "504.29"); a leading apostrophe forces Excel to keep them as
# literal text instead of auto-converting them to floating point numbers
# (which would corrupt the exact displayed text, e.g. trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '54.422.22'
$ws.Range('E2').Formula = '  +0.75%  '
$ws.Range('D3').Formula = '2.285.93'
$ws.Range('E3').Formula = '  +1.74%  '
$ws.Range('E4').Formula = '  -0.84%  '
$ws.Range('D5').Formula = '''504.29'
$ws.Range('E5').Formula = '  +1.71%  '
$ws.Range('D6').Formula = '''129.98'
$ws.Range('E6').Formula = '  +2.06%  '
$ws.Range('E7').Formula = '  -0.17%  '
$ws.Range('E8').Formula = '  +0.25%  '
$ws.Range('D9').Formula = '''0.0959'
$ws.Range('E9').Formula = '  +1.96%  '
$ws.Range('E10').Formula = '  +0.96%  '
$ws.Range('D11').Formula = '''0.336'
$ws.Range('E11').Formula = '  +4.45%  '
$ws.Range('D12').Formula = '''4.73'
$ws.Range('E12').Formula = '  +2.20%  '
$ws.Range('D13').Formula = '2.691.67'
$ws.Range('E13').Formula = '  +0.82%  '
$ws.Range('D14').Formula = '''22.98'
$ws.Range('E14').Formula = '  +6.19%  '
$ws.Range('D15').Formula = '54.374.83'
$ws.Range('E15').Formula = '  +0.63%  '
$ws.Range('E16').Formula = '  +0.45%  '
$ws.Range('D17').Formula = '2.280.09'
$ws.Range('E17').Formula = '  -0.19%  '
$ws.Range('D18').Formula = '''10.32'
$ws.Range('E18').Formula = '  +4.14%  '
$ws.Range('D19').Formula = '''4.15'
$ws.Range('E19').Formula = '  +2.49%  '
$ws.Range('D20').Formula = '''305.70'
$ws.Range('E20').Formula = '  +0.69%  '
$ws.Range('E21').Formula = '  -0.33%  '
$ws.Range('E22').Formula = '  +0.22%  '
$ws.Range('E23').Formula = '  -3.07%  '
$ws.Range('D24').Formula = '''0.999'
$ws.Range('E24').Formula = '  -0.48%  '
$ws.Range('D25').Formula = '''0.152'
$ws.Range('E25').Formula = '  +1.99%  '
$ws.Range('D26').Formula = '''7.35'
$ws.Range('E26').Formula = '  +2.97%  '
$ws.Range('D27').Formula = '''173.51'
$ws.Range('E27').Formula = '  +4.90%  '
$ws.Range('E28').Formula = '  +1.73%  '
$ws.Range('E29').Formula = '  +2.29%  '
$ws.Range('D30').Formula = '0.0₃0695'
$ws.Range('E30').Formula = '  +1.82%  '
$ws.Range('D31').Formula = '''1.10'
$ws.Range('E31').Formula = '  +2.62%  '
$ws.Range('E32').Formula = '  -0.02%  '
$ws.Range('D33').Formula = '''17.84'
$ws.Range('E33').Formula = '  +1.14%  '
$ws.Range('D34').Formula = '''0.982'
$ws.Range('E34').Formula = '  +12.60%  '
$ws.Range('D35').Formula = '''0.996'
$ws.Range('E35').Formula = '  +0.62%  '
$ws.Range('D36').Formula = '''1.21'
$ws.Range('E36').Formula = '  +1.64%  '
$ws.Range('E37').Formula = '  +3.61%  '
$ws.Range('D38').Formula = '''0.375'
$ws.Range('E38').Formula = '  -0.32%  '
$ws.Range('D39').Formula = '''1.42'
$ws.Range('E39').Formula = '  +1.31%  '
$ws.Range('E40').Formula = '  +1.37%  '
$ws.Range('D41').Formula = '''4.87'
$ws.Range('E41').Formula = '  +1.22%  '
$ws.Range('D42').Formula = '''125.02'
$ws.Range('E42').Formula = '  -0.91%  '
$ws.Range('D43').Formula = '''0.0498'
$ws.Range('E43').Formula = '  +3.58%  '
$ws.Range('E44').Formula = '  +0.60%  '
$ws.Range('B45').Formula = 'Bittensor'
$ws.Range('C45').Formula = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Formula = '''244.44'
$ws.Range('E45').Formula = '  +2.84%  '
$ws.Range('B46').Formula = 'Mantle'
$ws.Range('C46').Formula = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Formula = '''0.550'
$ws.Range('E46').Formula = '  +0.70%  '
$ws.Range('E47').Formula = '  -0.33%  '
$ws.Range('E48').Formula = '  +1.28%  '
$ws.Range('E49').Formula = '  +0.96%  '
$ws.Range('D50').Formula = '''16.55'
$ws.Range('E50').Formula = '  +1.50%  '
$ws.Range('E51').Formula = '  +0.05%  '
